$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-6 with new content
$ws.Range("A2").Value = "武汉"
$ws.Range("B2").Value = "武汉[chinese]"

$ws.Range("A3").Value = "name_1718273788838"
$ws.Range("B3").Value = "{{name}}"

$ws.Range("A4").Value = "address_1718273788839"
$ws.Range("B4").Value = "{{address}}"

$ws.Range("A5").Value = "我的名字："
$ws.Range("B5").Value = "我的名字：[chinese]"

$ws.Range("A6").Value = "东湖"
$ws.Range("B6").Value = "东湖[chinese]"

# Remove old rows 7-9 (no longer present in the sheet)
$ws.Range("A7:B9").ClearContents()
